# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    3  = 3174
    4  = 231
    5  = 125
    7  = 1678
    8  = 1627
    9  = 464
    10 = 368
    16 = 238
    19 = 24
    20 = 17
    21 = 51
    22 = 21
    23 = 376
    24 = 202
    25 = 101
    26 = 30
    27 = 8
    28 = 22
    29 = 209
    30 = 2146
    33 = 466
    34 = 320
    36 = 425
    38 = 343
    40 = 513
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
